$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lowercase the UUID-like "ID" values in column A (rows 2-8)
for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $orig = [string]$cell.Value2
    $cell.Value = $orig.ToLower()
}

# Apply new font/alignment style to A2:A8
$rng = $ws.Range("A2:A8")
$rng.Font.Name = "Arial Unicode MS"
$rng.Font.Size = 10
$rng.Font.Color = 0
$rng.VerticalAlignment = -4108  # xlCenter

# Autofit column B (Name) which contains long padded strings
$ws.Columns.Item(2).AutoFit() | Out-Null

# Move selection to A8
$ws.Range("A8").Select() | Out-Null
